$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "z"
$ws.Range("C1").Value = "y"
$ws.Range("D1").Value = "x"
$ws.Range("E1").Value = "mass"

$ws.Range("A1").Select()
